$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the column B (percentage) values in place. Row numbers below
#     refer to the *current* row positions, before the two rows are
#     deleted at the end of this script. ---

$ws.Range("B1").Value = 99.89626556016597
$ws.Range("B2").Value = 90.55499268265226
$ws.Range("B3").Value = 92.18697829716194
$ws.Range("B4").Value = 0.12558306422676
$ws.Range("B5").Value = 99.93240959783711
$ws.Range("B6").Value = 99.86126526082131
$ws.Range("B7").Value = 0.03562945368171
$ws.Range("B8").Value = 99.84097535117944
$ws.Range("B9").Value = 62.14399731160212
$ws.Range("B10").Value = 0.05029337803855
$ws.Range("B11").Value = 0.16602520564485
$ws.Range("B12").Value = 0.14882797966017
$ws.Range("B13").Value = 99.71719457013576
$ws.Range("B14").Value = 0.0204081632653
$ws.Range("B15").Value = 74.8594731621294
$ws.Range("B16").Value = 99.90435198469632
$ws.Range("B17").Value = 95.27236149798824
$ws.Range("B18").Value = 0.14529124290054
$ws.Range("B19").Value = 5.47242411286874
# Row 20 (GGAAGACCTGATACC) is removed entirely below, so no value update here.
$ws.Range("B21").Value = 0.11597564511452
$ws.Range("B22").Value = 90.09179521150752
$ws.Range("B23").Value = 99.84799864887688
$ws.Range("B24").Value = 0.08751969193068
$ws.Range("B25").Value = 0.10712372790573
$ws.Range("B26").Value = 99.43310657596372
$ws.Range("B27").Value = 79.65599051008304
$ws.Range("B28").Value = 0.11778563015312
# Row 29 (TCTTTAATCCAGATA) is removed entirely below, so no value update here.
$ws.Range("B30").Value = 0.02645502645502
$ws.Range("B31").Value = 0.0605815831987
$ws.Range("B32").Value = 52.58803283028859

# --- Drop the two rows removed from the dataset (GGAAGACCTGATACC / row 20,
#     TCTTTAATCCAGATA / row 29). Delete the higher row index first so the
#     other index stays valid after the shift. ---
$ws.Rows.Item(29).Delete()
$ws.Rows.Item(20).Delete()
